$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($ws, $addr, $val)
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-CellText $ws 'D2' '34.603.89'
Set-CellText $ws 'E2' '  +2.12%  '
Set-CellText $ws 'D3' '1.789.36'
Set-CellText $ws 'E4' '  +0.00%  '
Set-CellText $ws 'D5' '224.29'
Set-CellText $ws 'E5' '  -0.31%  '
Set-CellText $ws 'E6' '  +0.68%  '
Set-CellText $ws 'E7' '  +0.04%  '
Set-CellText $ws 'D8' '32.57'
Set-CellText $ws 'E8' '  +6.54%  '
Set-CellText $ws 'D9' '0.283'
Set-CellText $ws 'E9' '  +2.25%  '
Set-CellText $ws 'D10' '0.0668'
Set-CellText $ws 'E10' '  +1.25%  '
Set-CellText $ws 'D11' '0.0936'
Set-CellText $ws 'E11' '  +1.33%  '
Set-CellText $ws 'D12' '2.047.24'
Set-CellText $ws 'E12' '  +0.51%  '
Set-CellText $ws 'D13' '10.99'
Set-CellText $ws 'E13' '  +10.37%  '
Set-CellText $ws 'D14' '1.779.68'
Set-CellText $ws 'E14' '  -0.05%  '
Set-CellText $ws 'B15' 'WrappedBTC'
Set-CellText $ws 'C15' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-CellText $ws 'D15' '34.626.79'
Set-CellText $ws 'E15' '  +2.22%  '
Set-CellText $ws 'B16' 'Polygon'
Set-CellText $ws 'C16' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-CellText $ws 'D16' '0.632'
Set-CellText $ws 'E16' '  +0.88%  '
Set-CellText $ws 'D17' '4.28'
Set-CellText $ws 'E17' '  +2.29%  '
Set-CellText $ws 'D18' '68.67'
Set-CellText $ws 'E18' '  +0.45%  '
Set-CellText $ws 'D19' '253.74'
Set-CellText $ws 'E19' '  +1.15%  '
Set-CellText $ws 'D20' '0.0₃0765'
Set-CellText $ws 'E20' '  +3.72%  '
Set-CellText $ws 'E21' '  -0.01%  '
Set-CellText $ws 'D22' '10.39'
Set-CellText $ws 'E22' '  +0.96%  '
Set-CellText $ws 'D23' '4.23'
Set-CellText $ws 'E23' '  +0.44%  '
Set-CellText $ws 'E24' '  -1.06%  '
Set-CellText $ws 'D25' '159.46'
Set-CellText $ws 'E25' '  +0.26%  '
Set-CellText $ws 'D26' '16.36'
Set-CellText $ws 'E26' '  -0.70%  '
Set-CellText $ws 'D27' '7.07'
Set-CellText $ws 'E27' '  +2.07%  '
Set-CellText $ws 'E28' '  +0.35%  '
Set-CellText $ws 'E29' '  +0.10%  '
Set-CellText $ws 'B30' 'Hedera'
Set-CellText $ws 'C30' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-CellText $ws 'D30' '0.0516'
Set-CellText $ws 'E30' '  +0.61%  '
Set-CellText $ws 'B31' 'Filecoin'
Set-CellText $ws 'C31' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-CellText $ws 'D31' '3.75'
Set-CellText $ws 'E31' '  -1.10%  '
Set-CellText $ws 'E32' '  +0.15%  '
Set-CellText $ws 'E33' '  +0.75%  '
Set-CellText $ws 'E34' '  +2.92%  '
Set-CellText $ws 'D35' '1.441.62'
Set-CellText $ws 'E35' '  -2.66%  '
Set-CellText $ws 'E36' '  -0.19%  '
Set-CellText $ws 'E37' '  +2.39%  '
Set-CellText $ws 'D38' '0.628'
Set-CellText $ws 'E38' '  -0.51%  '
Set-CellText $ws 'D39' '82.96'
Set-CellText $ws 'E39' '  -0.38%  '
Set-CellText $ws 'D40' '2.80'
Set-CellText $ws 'E40' '  +3.98%  '
Set-CellText $ws 'E41' '  -0.26%  '
Set-CellText $ws 'D42' '0.898'
Set-CellText $ws 'E42' '  +1.40%  '
Set-CellText $ws 'E43' '  -0.36%  '
Set-CellText $ws 'E44' '  -0.73%  '
Set-CellText $ws 'D45' '5.89'
Set-CellText $ws 'E45' '  +2.08%  '
Set-CellText $ws 'E46' '  -1.65%  '
Set-CellText $ws 'D47' '1.942.38'
Set-CellText $ws 'E47' '  +0.28%  '
Set-CellText $ws 'D48' '12.01'
Set-CellText $ws 'E48' '  +0.36%  '
Set-CellText $ws 'E49' '  -0.02%  '
Set-CellText $ws 'D50' '103.20'
Set-CellText $ws 'E50' '  +5.77%  '
Set-CellText $ws 'E51' '  +4.57%  '
